$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("React-MainConcepts")
$nl = [char]10

# ---------------------------------------------------------------------------
# Column B width tweak (fits new "Components and Props" / "State and Lifecycle")
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 21.85546875

# ---------------------------------------------------------------------------
# Row 12: "Rendering Elements" -> "React Elements", new Comments cell
# ---------------------------------------------------------------------------
$b12 = $ws.Range("B12")
$b12.Value = "React Elements"
$b12.Characters(1,5).Font.Bold = $true
$b12.Characters(7,8).Font.Bold = $true

$ws.Range("C12").Value = "Smallest building block, plain object, cheap to create"

$d12 = $ws.Range("D12")
$d12.Value = "ReactDOM.render(${nl}  element,${nl}  document.getElementById('root')${nl});"
$d12.WrapText = $true
$d12.VerticalAlignment = -4108
$d12.Font.Name = "Consolas"

$ws.Rows.Item(12).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 13-18
# ---------------------------------------------------------------------------
$ws.Range("C13").Value = "Describes what you want to see on screen"
$ws.Range("C14").Value = "React elements are immutable."
$ws.Range("C15").Value = "React DOM takes care of updating the DOM to match the React elements."
$ws.Range("C16").Value = "The only way to update the UI is to create a new element, and pass it to ReactDOM.render()."

$c17 = $ws.Range("C17")
$c17.Value = "In practice, most React apps only call ReactDOM.render() once.  ==>> code gets encapsulated into stateful components."
$c17.Characters(39,18).Font.Bold = $true
$c17.WrapText = $true
$ws.Rows.Item(17).RowHeight = 30

$ws.Range("C18").Value = "React Only Updates What" + [char]8217 + "s Necessary"

# ---------------------------------------------------------------------------
# Row 20-28 : Components and Props
# ---------------------------------------------------------------------------
$b20 = $ws.Range("B20")
$b20.Value = "Components and Props"
$b20.Font.Bold = $true

$ws.Range("C20").Value = "Components let you split the UI into independent, reusable pieces, and think about each piece in isolation."

$c21 = $ws.Range("C21")
$c21.Value = "Conceptually, components are like JavaScript functions, I/P(Props)==> O/P(React elements describing what should appear on the screen.)"
$c21.Characters(57,10).Font.Bold = $true
$c21.Characters(71,3).Font.Bold = $true

$ws.Range("C22").Value = "Function component"
$d22 = $ws.Range("D22")
$d22.WrapText = $true

$ws.Range("C23").Value = "class component"
$d23 = $ws.Range("D23")
$d23.Value = "class Welcome extends React.Component {${nl}  render() {${nl}    return <h1>Hello, {this.props.name}</h1>;${nl}  }${nl}}"
$d23.WrapText = $true
$ws.Rows.Item(23).RowHeight = 75

$ws.Range("C24").Value = "can render DOM tags + user defined components"

$c25 = $ws.Range("C25")
$c25.Value = "1.  We call ReactDOM.render() with the <Welcome name=`"Sara`" /> element.${nl}2.  React calls the Welcome component with {name: 'Sara'} as the props.${nl}3.  Our Welcome component returns a <h1>Hello, Sara</h1> element as the result.${nl}4.  React DOM efficiently updates the DOM to match <h1>Hello, Sara</h1>.${nl}${nl}Note: Always start component names with a capital letter."
$c25.Characters(298,5).Font.Bold = $true
$c25.WrapText = $true
$c25.VerticalAlignment = -4160

$d25 = $ws.Range("D25")
$d25.Value = "function Welcome(props) {${nl}  return <h1>Hello, {props.name}</h1>;${nl}}${nl}const element = <Welcome name=`"Sara`" />;${nl}ReactDOM.render(${nl}  element,${nl}  document.getElementById('root')${nl});"
$d25.WrapText = $true
$ws.Rows.Item(25).RowHeight = 135

$ws.Range("C26").Value = "Components can refer to other components in their output."
$ws.Range("C27").Value = "Props naming ==>> component" + [char]8217 + "s own point of view rather than the context in which it is being used."
$ws.Range("C28").Value = "Props are Read-Only"

# ---------------------------------------------------------------------------
# Row 30-36 : State and Lifecycle
# ---------------------------------------------------------------------------
$b30 = $ws.Range("B30")
$b30.Value = "State and Lifecycle"
$b30.Font.Bold = $true

$ws.Range("C30").Value = "State is similar to props, but it is private and fully controlled by the component."
$d30 = $ws.Range("D30")
$d30.Value = "// Wrong${nl}this.state.comment = 'Hello';${nl}// Correct${nl}this.setState({comment: 'Hello'});"
$d30.WrapText = $true
$d30.VerticalAlignment = -4160
# Multi-line wrapped content would otherwise force an auto row height; row 30
# stays default height in the source workbook (content lives in a merged
# D30:D35 block), so re-fit it back down.
$ws.Rows.Item(30).AutoFit()

$ws.Range("C31").Value = "Do Not Modify State Directly"
$d31 = $ws.Range("D31")
$d31.WrapText = $true
$d31.VerticalAlignment = -4160
$ws.Rows.Item(31).RowHeight = 15

$ws.Range("C32").Value = "The only place where you can assign this.state is the constructor."
$d32 = $ws.Range("D32")
$d32.WrapText = $true
$d32.VerticalAlignment = -4160

$ws.Range("C33").Value = "State Updates are Merged"
$d33 = $ws.Range("D33")
$d33.WrapText = $true
$d33.VerticalAlignment = -4160

$ws.Range("C34").Value = "The merging is shallow, so this.setState({comments}) leaves this.state.posts intact, but completely replaces this.state.comments."
$d34 = $ws.Range("D34")
$d34.WrapText = $true
$d34.VerticalAlignment = -4160

$ws.Range("C35").Value = "A component may choose to pass its state down as props to its child components:"
$d35 = $ws.Range("D35")
$d35.WrapText = $true
$d35.VerticalAlignment = -4160

$ws.Range("C36").Value = "You can use stateless components inside stateful components, and vice versa."

# ---------------------------------------------------------------------------
# Merge the D30:D35 comment block (mirrors existing D2:D5 merge)
# ---------------------------------------------------------------------------
$ws.Range("D30:D35").Merge()

# ---------------------------------------------------------------------------
# Selection / view state -> land on the newly-added last row
# ---------------------------------------------------------------------------
$ws.Range("C36").Select()
$excel.ActiveWindow.ScrollRow = 24
